$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7205305911032269
$ws.Range("D2").Value = 0.01126933737403135
$ws.Range("E2").Value = 0.0366412189657721
$ws.Range("F2").Value = 4.528730623816841
$ws.Range("G2").Value = 0.002641910959350843
$ws.Range("I2").Value = 2.494931622134104
$ws.Range("J2").Value = 0.1014389295986393
$ws.Range("K2").Value = 1.468597778396258
$ws.Range("L2").Value = 0.5926008526189719

$ws.Range("B3").Value = 0.7133901413602644
$ws.Range("D3").Value = 0.009876196763748624
$ws.Range("E3").Value = 0.03610054927880046
$ws.Range("F3").Value = 4.480448900386264
$ws.Range("G3").Value = 0.002647047282163658
$ws.Range("I3").Value = 2.475614974372661
$ws.Range("J3").Value = 0.1009483339147441
$ws.Range("K3").Value = 1.419528241681235
$ws.Range("L3").Value = 0.5836601267851336

$ws.Range("B4").Value = 0.7095018409301019
$ws.Range("D4").Value = 0.009018767946329831
$ws.Range("E4").Value = 0.03576334569520778
$ws.Range("F4").Value = 4.452224354894994
$ws.Range("G4").Value = 0.002650367787376107
$ws.Range("I4").Value = 2.464290087476286
$ws.Range("J4").Value = 0.1006361114359202
$ws.Range("K4").Value = 1.390476022379204
$ws.Range("L4").Value = 0.5785223269987512

$ws.Range("B5").Value = 0.7080423126041921
$ws.Range("D5").Value = 0.008668795280314612
$ws.Range("E5").Value = 0.0356245946535374
$ws.Range("F5").Value = 4.441079086594129
$ws.Range("G5").Value = 0.002651762999087996
$ws.Range("I5").Value = 2.459809425440454
$ws.Range("J5").Value = 0.1005060896942847
$ws.Range("K5").Value = 1.378907239804505
$ws.Range("L5").Value = 0.576517031528212

$ws.Range("B6").Value = 0.7078075159496251
$ws.Range("D6").Value = 0.008610647042551989
$ws.Range("E6").Value = 0.03560147368720923
$ws.Range("F6").Value = 4.439249930639988
$ws.Range("G6").Value = 0.002651997218637968
$ws.Range("I6").Value = 2.459073517034625
$ws.Range("J6").Value = 0.1004843305667817
$ws.Range("K6").Value = 1.377002559291327
$ws.Range("L6").Value = 0.5761893908095459

$ws.Range("B7").Value = 0.7094816507976986
$ws.Range("D7").Value = 0.009014050441479071
$ws.Range("E7").Value = 0.03576147989627465
$ws.Range("F7").Value = 4.452072603640772
$ws.Range("G7").Value = 0.002650386432990742
$ws.Range("I7").Value = 2.464229116334607
$ws.Range("J7").Value = 0.1006343692387492
$ws.Range("K7").Value = 1.390318908231649
$ws.Range("L7").Value = 0.5784949250517855

$ws.Range("B8").Value = 0.7179657829609027
$ws.Range("D8").Value = 0.01078937326585816
$ws.Range("E8").Value = 0.03645586778909138
$ws.Range("F8").Value = 4.511787750192894
$ws.Range("G8").Value = 0.002643647432236025
$ws.Range("I8").Value = 2.488159848667991
$ws.Range("J8").Value = 0.1012720397401807
$ws.Range("K8").Value = 1.45145480876721
$ws.Range("L8").Value = 0.58944501531208

$ws.Range("B9").Value = 0.7385287902021105
$ws.Range("D9").Value = 0.0142571085679748
$ws.Range("E9").Value = 0.03777713139314898
$ws.Range("F9").Value = 4.640205006299709
$ws.Range("G9").Value = 0.002631749213273808
$ws.Range("I9").Value = 2.53935770432517
$ws.Range("J9").Value = 0.1024363682826639
$ws.Range("K9").Value = 1.579917263000766
$ws.Range("L9").Value = 0.6137153226802923

$ws.Range("B10").Value = 0.7560201282798857
$ws.Range("D10").Value = 0.01680015808717883
$ws.Range("E10").Value = 0.03872476847989681
$ws.Range("F10").Value = 4.741526312787613
$ws.Range("G10").Value = 0.002623801463189634
$ws.Range("I10").Value = 2.579608420606988
$ws.Range("J10").Value = 0.1032408368470685
$ws.Range("K10").Value = 1.679584427791781
$ws.Range("L10").Value = 0.6332627396915598

$ws.Range("B11").Value = 0.7644932929803758
$ws.Range("D11").Value = 0.01795683597941178
$ws.Range("E11").Value = 0.03915119937895639
$ws.Range("F11").Value = 4.789151219041855
$ws.Range("G11").Value = 0.002620356294351905
$ws.Range("I11").Value = 2.598499198694043
$ws.Range("J11").Value = 0.103596081208293
$ws.Range("K11").Value = 1.726087035905152
$ws.Range("L11").Value = 0.6425305780236101

$ws.Range("B12").Value = 0.7677759030159166
$ws.Range("D12").Value = 0.01839487611705692
$ws.Range("E12").Value = 0.03931203495608315
$ws.Range("F12").Value = 4.80740714492336
$ws.Range("G12").Value = 0.002619076041969759
$ws.Range("I12").Value = 2.605736616519451
$ws.Range("J12").Value = 0.1037290892305993
$ws.Range("K12").Value = 1.743864549279181
$ws.Range("L12").Value = 0.646094240032852

$ws.Range("B13").Value = 0.7670656463876924
$ws.Range("D13").Value = 0.018300534055399
$ws.Range("E13").Value = 0.03927742440739479
$ws.Range("F13").Value = 4.803465542764656
$ws.Range("G13").Value = 0.002619350685974713
$ws.Range("I13").Value = 2.604174171478334
$ws.Range("J13").Value = 0.1037005105652815
$ws.Range("K13").Value = 1.740028361571433
$ws.Range("L13").Value = 0.6453243321352033

$ws.Range("B14").Value = 0.764761873392473
$ws.Range("D14").Value = 0.01799287285805207
$ws.Range("E14").Value = 0.03916444419405885
$ws.Range("F14").Value = 4.790648702427774
$ws.Range("G14").Value = 0.002620250479940365
$ws.Range("I14").Value = 2.599092941766514
$ws.Range("J14").Value = 0.1036070540827936
$ws.Range("K14").Value = 1.727546232030022
$ws.Range("L14").Value = 0.6428226770393906

$ws.Range("B15").Value = 0.7633603769486967
$ws.Range("D15").Value = 0.01780442722526487
$ws.Range("E15").Value = 0.03909515734925861
$ws.Range("F15").Value = 4.782826878816934
$ws.Range("G15").Value = 0.002620804797888479
$ws.Range("I15").Value = 2.595991479072879
$ws.Range("J15").Value = 0.1035496126998297
$ws.Range("K15").Value = 1.719922464208025
$ws.Range("L15").Value = 0.6412973943675695

$ws.Range("B16").Value = 0.7554767556780178
$ws.Range("D16").Value = 0.01672456781919607
$ws.Range("E16").Value = 0.03869680867185643
$ws.Range("F16").Value = 4.738444845785153
$ws.Range("G16").Value = 0.002624030028905282
$ws.Range("I16").Value = 2.578385580587209
$ws.Range("J16").Value = 0.103217407374844
$ws.Range("K16").Value = 1.676568842910143
$ws.Range("L16").Value = 0.6326646331556844

$ws.Range("B17").Value = 0.750772470763394
$ws.Range("D17").Value = 0.01606210536301944
$ws.Range("E17").Value = 0.03845126014517675
$ws.Range("F17").Value = 4.711611277700968
$ws.Range("G17").Value = 0.002626052128730333
$ws.Range("I17").Value = 2.567733910434853
$ws.Range("J17").Value = 0.1030108842795272
$ws.Range("K17").Value = 1.650271276493868
$ws.Range("L17").Value = 0.6274650021566544

$ws.Range("B18").Value = 0.7481152975458656
$ws.Range("D18").Value = 0.01568105531747932
$ws.Range("E18").Value = 0.03830958758070047
$ws.Range("F18").Value = 4.696321588420204
$ws.Range("G18").Value = 0.002627231224426875
$ws.Range("I18").Value = 2.561661968561097
$ws.Range("J18").Value = 0.1028910878389677
$ws.Range("K18").Value = 1.635255085825833
$ws.Range("L18").Value = 0.6245096596551321

$ws.Range("B19").Value = 0.7472239798217402
$ws.Range("D19").Value = 0.015552033343738
$ws.Range("E19").Value = 0.03826154344484678
$ws.Range("F19").Value = 4.691169513025272
$ws.Range("G19").Value = 0.002627633204503433
$ws.Range("I19").Value = 2.559615481407121
$ws.Range("J19").Value = 0.1028503526238076
$ws.Range("K19").Value = 1.630189643784632
$ws.Range("L19").Value = 0.6235150989263616

$ws.Range("B20").Value = 0.7512682209894308
$ws.Range("D20").Value = 0.01613262720435671
$ws.Range("E20").Value = 0.03847744451742408
$ws.Range("F20").Value = 4.71445281896527
$ws.Range("G20").Value = 0.002625835213875427
$ws.Range("I20").Value = 2.568862143303718
$ws.Range("J20").Value = 0.1030329733447717
$ws.Range("K20").Value = 1.653059362398693
$ws.Range("L20").Value = 0.628014852557385

$ws.Range("B21").Value = 0.7654365406727663
$ws.Range("D21").Value = 0.01808323911537002
$ws.Range("E21").Value = 0.03919764651617541
$ws.Range("F21").Value = 4.794407302823004
$ws.Range("G21").Value = 0.002619985528496771
$ws.Range("I21").Value = 2.600583142057516
$ws.Range("J21").Value = 0.1036345454405954
$ws.Range("K21").Value = 1.731207971514607
$ws.Range("L21").Value = 0.6435560036591426

$ws.Range("B22").Value = 0.7751276193464207
$ws.Range("D22").Value = 0.01935828536821305
$ws.Range("E22").Value = 0.03966459488380991
$ws.Range("F22").Value = 4.847953146221272
$ws.Range("G22").Value = 0.002616304339057569
$ws.Range("I22").Value = 2.621803817222556
$ws.Range("J22").Value = 0.1040188870039884
$ws.Range("K22").Value = 1.783262077676852
$ws.Range("L22").Value = 0.6540286305620953

$ws.Range("B23").Value = 0.7699159228716894
$ws.Range("D23").Value = 0.01867773166144104
$ws.Range("E23").Value = 0.03941571021547219
$ws.Range("F23").Value = 4.819256292528138
$ws.Range("G23").Value = 0.002618256116664714
$ws.Range("I23").Value = 2.61043304794579
$ws.Range("J23").Value = 0.1038145554114669
$ws.Range("K23").Value = 1.755389971757666
$ws.Range("L23").Value = 0.6484102756098196

$ws.Range("B24").Value = 0.7510439447032979
$ws.Range("D24").Value = 0.01610074488188928
$ws.Range("E24").Value = 0.03846560813324551
$ws.Range("F24").Value = 4.713167730582512
$ws.Range("G24").Value = 0.002625933229328406
$ws.Range("I24").Value = 2.56835190775405
$ws.Range("J24").Value = 0.1030229901936837
$ws.Range("K24").Value = 1.651798549133304
$ws.Range("L24").Value = 0.6277661593120314

$ws.Range("B25").Value = 0.7325466862129133
$ws.Range("D25").Value = 0.01332008284340702
$ws.Range("E25").Value = 0.03742386336606085
$ws.Range("F25").Value = 4.604245608411446
$ws.Range("G25").Value = 0.00263482793897698
$ws.Range("I25").Value = 2.525046849273252
$ws.Range("J25").Value = 0.1021304489230799
$ws.Range("K25").Value = 1.544240801762783
$ws.Range("L25").Value = 0.606849006363305

Write-Output "Applied 380 kV case values"